# Auto-generated: apply cell text updates for cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.023.11"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "2.476.60"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'583.78"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").Value = "'169.33"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.95%  "
$ws.Range("D9").Value = "2.474.87"
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("D10").Value = "'0.136"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("E13").Value = "  -3.54%  "
$ws.Range("D14").Value = "2.921.07"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("E15").Value = "  -2.99%  "
$ws.Range("D16").Value = "66.939.59"
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("E17").Value = "  -3.15%  "
$ws.Range("D18").Value = "2.491.38"
$ws.Range("E18").Value = "  -2.19%  "
$ws.Range("D19").Value = "'10.96"
$ws.Range("E19").Value = "  -6.10%  "
$ws.Range("D20").Value = "'7.37"
$ws.Range("E20").Value = "  -8.28%  "
$ws.Range("D21").Value = "'349.51"
$ws.Range("E21").Value = "  -3.87%  "
$ws.Range("D22").Value = "'3.98"
$ws.Range("E22").Value = "  -3.95%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").Value = "'68.40"
$ws.Range("E24").Value = "  -4.69%  "
$ws.Range("E25").Value = "  -6.63%  "
$ws.Range("E26").Value = "  -2.64%  "
$ws.Range("D27").Value = "'9.11"
$ws.Range("E27").Value = "  -6.17%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -35.87%  "
$ws.Range("D29").Value = "2.595.40"
$ws.Range("E29").Value = "  -2.51%  "
$ws.Range("D30").Value = "0.0₃0895"
$ws.Range("E30").Value = "  -4.69%  "
$ws.Range("D31").Value = "'505.88"
$ws.Range("E31").Value = "  -4.47%  "
$ws.Range("D32").Value = "'7.61"
$ws.Range("E32").Value = "  -7.31%  "
$ws.Range("E33").Value = "  -4.65%  "
$ws.Range("E34").Value = "  -4.50%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'158.50"
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("E37").Value = "  -9.67%  "
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("D39").Value = "'18.16"
$ws.Range("E39").Value = "  -6.14%  "
$ws.Range("E40").Value = "  -7.13%  "
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("E42").Value = "  -4.80%  "
$ws.Range("D43").Value = "'4.78"
$ws.Range("E43").Value = "  -4.71%  "
$ws.Range("E44").Value = "  -4.83%  "
$ws.Range("D45").Value = "'2.34"
$ws.Range("E45").Value = "  -3.46%  "
$ws.Range("D46").Value = "'38.81"
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").Value = "'141.05"
$ws.Range("E47").Value = "  -4.30%  "
$ws.Range("E48").Value = "  -6.95%  "
$ws.Range("D49").Value = "'0.510"
$ws.Range("E49").Value = "  -6.35%  "
$ws.Range("D50").Value = "0.0₆0249"
$ws.Range("E50").Value = "  -7.90%  "
$ws.Range("D51").Value = "'0.0725"
$ws.Range("E51").Value = "  -1.94%  "
